# Phemex futures trades workbook fix:
#  - rename/shift header columns (insert "Last Price", rename tier/possible-trade
#    columns to Move / Move percentage / Tier, rename SL80%/CreateDate to
#    On exchange / Note)
#  - rename ticker BTCUSDPERP -> uBTCUSD
#  - replace sample trade numbers and drop the helper TRUE()/FALSE() tier
#    formulas in favour of blank (but still boolean-formatted) cells
#  - On exchange column now holds the literal text "A1" instead of a number

$wb = $excel.ActiveWorkbook

function Set-TradeSheetHeaders($ws) {
    $ws.Range("B1").Value2 = "Last Price"
    $ws.Range("C1").Value2 = "Entry price"
    $ws.Range("D1").Value2 = "Stop loss"
    $ws.Range("E1").Value2 = "Profit target 1"
    $ws.Range("F1").Value2 = "Position"
    $ws.Range("G1").Value2 = "Leverage"
    $ws.Range("H1").Value2 = "Move"
    $ws.Range("I1").Value2 = "Move percentage"
    $ws.Range("J1").Value2 = "Tier"
    $ws.Range("K1").Value2 = "On exchange"
    $ws.Range("L1").Value2 = "Note"
}

function Set-TradeSheetRow2($ws, $entry, $stop, $target) {
    $ws.Range("A2").Value2 = "uBTCUSD"
    $ws.Range("B2").ClearContents()
    $ws.Range("C2").Value2 = $entry
    $ws.Range("D2").Value2 = $stop
    $ws.Range("E2").Value2 = $target
    $ws.Range("F2").Value2 = 0.001
    $ws.Range("G2").Value2 = 20
    $ws.Range("G2").NumberFormat = "General"
    $ws.Range("H2").ClearContents()
    $ws.Range("I2").ClearContents()
    $ws.Range("J2").ClearContents()
    $ws.Range("K2").Value2 = "A1"
}

# --- TradesLong ---
$ws1 = $wb.Worksheets.Item("TradesLong")
Set-TradeSheetHeaders $ws1
$ws1.Range("M1").Clear()
$ws1.Range("N1").Clear()
Set-TradeSheetRow2 $ws1 15000 14900 17200
$ws1.Range("M2").Clear()
$ws1.Range("D5").Select()

# --- TradesShort ---
$ws2 = $wb.Worksheets.Item("TradesShort")
Set-TradeSheetHeaders $ws2
$ws2.Range("M1").ClearContents()
$ws2.Range("N1").ClearContents()
Set-TradeSheetRow2 $ws2 70000 70100 68000
$ws2.Range("M2").Clear()
$ws2.Range("O4").Select()

# TradesLong is the tab that should remain active/selected.
$ws1.Activate()
